$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'269.65"
$ws.Range("E2").Value = "'3.03%"
$ws.Range("E3").Value = "'-1.71%"
$ws.Range("D4").Value = "'4.712"
$ws.Range("E4").Value = "'0.11%"
$ws.Range("D5").Value = "'0.06099"
$ws.Range("E5").Value = "'-1.78%"
$ws.Range("D6").Value = "'6.746"
$ws.Range("E6").Value = "'0.42%"
$ws.Range("D7").Value = "'0.8571"
$ws.Range("E7").Value = "'0.83%"
$ws.Range("D8").Value = "'0.8895"
$ws.Range("E8").Value = "'-2.77%"
$ws.Range("D9").Value = "'0.1430"
$ws.Range("E9").Value = "'1.35%"
$ws.Range("D10").Value = "'0.04968"
$ws.Range("E10").Value = "'6.44%"
$ws.Range("D11").Value = "'0.07108"
$ws.Range("E11").Value = "'0.29%"
$ws.Range("D12").Value = "'0.03183"
$ws.Range("E12").Value = "'1.02%"
$ws.Range("D13").Value = "'0.09037"
$ws.Range("E13").Value = "'-0.20%"
$ws.Range("D14").Value = "'0.001538"
$ws.Range("E14").Value = "'0.31%"
$ws.Range("E15").Value = "'-1.35%"
$ws.Range("D16").Value = "'0.005953"
$ws.Range("E16").Value = "'-2.67%"
$ws.Range("D17").Value = "'3.463"
$ws.Range("E17").Value = "'-0.22%"
$ws.Range("E18").Value = "'0.04%"
$ws.Range("D19").Value = "'2.244"
$ws.Range("E19").Value = "'2.95%"
$ws.Range("D21").Value = "'0.1300"
$ws.Range("E21").Value = "'-0.80%"
$ws.Range("D22").Value = "'3.839"
$ws.Range("E22").Value = "'-5.94%"
$ws.Range("D23").Value = "'0.04251"
$ws.Range("E23").Value = "'0.15%"
$ws.Range("D24").Value = "'0.001177"
$ws.Range("E24").Value = "'-2.88%"
$ws.Range("D25").Value = "'0.004148"
$ws.Range("E25").Value = "'0.36%"
$ws.Range("D26").Value = "'0.0001200"
$ws.Range("E26").Value = "'-0.04%"
$ws.Range("E27").Value = "'5.03%"
$ws.Range("D40").Value = "'0.03952"
$ws.Range("E40").Value = "'1.14%"
$ws.Range("E41").Value = "'0.61%"
$ws.Range("E42").Value = "'1.42%"
$ws.Range("D43").Value = "'0.002037"
$ws.Range("E43").Value = "'-6.74%"
$ws.Range("D44").Value = "'0.01179"
$ws.Range("E44").Value = "'-15.31%"
$ws.Range("D45").Value = "'0.00005136"
$ws.Range("E45").Value = "'-0.71%"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'-0.01%"
$ws.Range("B47").Value = "CoinbaseStockToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D47").Value = "'0.02448"
$ws.Range("E47").Value = "'-31.82%"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "'1.011"
$ws.Range("E48").Value = "'506.92%"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'-0.01%"
$ws.Range("E50").Value = "'-0.01%"
